# This script applies a data re-sync for the "Artfynd" sheet:
# row 3 <-> row 4 swap their full contents, and rows 5, 6, 8 rotate
# (new row5 = old row6, new row6 = old row8, new row8 = old row5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 110282846
$ws.Range("B3").Value = 103288
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 221144
$ws.Range("F3").Value = 'Grönpyrola'
$ws.Range("G3").Value = 'Pyrola chlorantha'
$ws.Range("H3").Value = 'Sw.'
$ws.Range("I3").Value = "'" + '10'
$ws.Range("J3").Value = 'plantor/tuvor'
$ws.Range("K3").Value = 'blomning'
$ws.Range("Q3").Value = 600839.9318167433
$ws.Range("R3").Value = 6613983.990819811

# --- Row 4 ---
$ws.Range("A4").Value = 110282820
$ws.Range("B4").Value = 89425
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 5442
$ws.Range("F4").Value = 'Tallticka'
$ws.Range("G4").Value = 'Porodaedalea pini'
$ws.Range("H4").Value = '(Brot.) Murrill'
$ws.Range("I4").Value = "'" + '1'
$ws.Range("J4").Value = 'fruktkroppar'
$ws.Range("K4").Value = ""
$ws.Range("Q4").Value = 600724.7123983201
$ws.Range("R4").Value = 6614086.574870056

# --- Row 5 ---
$ws.Range("A5").Value = 110282764
$ws.Range("B5").Value = 96348
$ws.Range("D5").Value = 'VU'
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = 'Knärot'
$ws.Range("G5").Value = 'Goodyera repens'
$ws.Range("H5").Value = '(L.) R. Br.'
$ws.Range("I5").Value = "'" + '10'
$ws.Range("J5").Value = 'plantor/tuvor'
$ws.Range("K5").Value = 'fullt utvecklade blad'
$ws.Range("Q5").Value = 600749.0751519018
$ws.Range("R5").Value = 6613971.934424319

# --- Row 6 ---
$ws.Range("A6").Value = 110282836
$ws.Range("B6").Value = 89793
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 4217
$ws.Range("F6").Value = 'Blodticka'
$ws.Range("G6").Value = 'Meruliopsis taxicola'
$ws.Range("H6").Value = '(Pers.:Fr.) Bondartsev'
$ws.Range("I6").Value = "'" + '1'
$ws.Range("J6").Value = 'mycel'
$ws.Range("K6").Value = ""
$ws.Range("Q6").Value = 600805.3583702671
$ws.Range("R6").Value = 6613969.910894822

# --- Row 8 ---
$ws.Range("A8").Value = 110282835
$ws.Range("B8").Value = 89405
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = 'Ullticka'
$ws.Range("G8").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H8").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I8").Value = "'" + '1'
$ws.Range("J8").Value = 'mycel'
$ws.Range("K8").Value = ""
$ws.Range("Q8").Value = 600805.3583702671
$ws.Range("R8").Value = 6613969.910894822
